$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct teacher name typos ("Adriano Nakamura" -> "Adriano")
$ws.Range("E5").Value = "Algoritmos Avançados / Adriano / 6"
$ws.Range("E7").Value = "Redes de computadores / Adriano / 4"

# Update teacher count
$ws.Range("A6").Value = "Número de professores: 16"
